$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and week-covering dates) ---
$ws.Range("A8").Value = "Volume 29   Number  51"
$ws.Range("C9").Value = "Report Covering the Week  12/19/2022  Through  12/25/2022"

# --- Cells changing from numeric to text dash markers ---
$ws.Range("C14").NumberFormat = "General"
$ws.Range("C17").NumberFormat = "General"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").NumberFormat = "General"
$ws.Range("C27").NumberFormat = "General"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").NumberFormat = "General"

# --- Cells changing from text dash markers back to numeric ---
$ws.Range("D30").NumberFormat = "#,##0"
$ws.Range("G30").NumberFormat = "#,##0"
$ws.Range("E30").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("H30").NumberFormat = "#,##0.0;""-""#,##0.0"

# --- Set cell values ---
$ws.Range("C14").Value = "0"
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -33.333333333333
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 14
$ws.Range("H16").Value = -64.285714285714
$ws.Range("I16").Value = 90
$ws.Range("J16").Value = 94
$ws.Range("K16").Value = -4.255319148936
$ws.Range("L16").Value = 40.625
$ws.Range("M16").Value = 60.714285714285
$ws.Range("N16").Value = -85.342019543973
$ws.Range("C17").Value = "0"
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = -100
$ws.Range("G17").Value = 20
$ws.Range("H17").Value = -60
$ws.Range("J17").Value = 126
$ws.Range("K17").Value = -10.31746031746
$ws.Range("L17").Value = 66.176470588235
$ws.Range("N17").Value = -29.813664596273
$ws.Range("C18").Value = 5
$ws.Range("E18").Value = 400
$ws.Range("F18").Value = 14
$ws.Range("G18").Value = 8
$ws.Range("H18").Value = 75
$ws.Range("I18").Value = 182
$ws.Range("J18").Value = 147
$ws.Range("K18").Value = 23.809523809523
$ws.Range("L18").Value = -3.191489361702
$ws.Range("M18").Value = 89.583333333333
$ws.Range("N18").Value = -84.971098265896
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = -45.454545454545
$ws.Range("F19").Value = 63
$ws.Range("G19").Value = 47
$ws.Range("H19").Value = 34.042553191489
$ws.Range("I19").Value = 660
$ws.Range("J19").Value = 539
$ws.Range("K19").Value = 22.448979591836
$ws.Range("L19").Value = 15.992970123022
$ws.Range("M19").Value = -10.326086956521
$ws.Range("N19").Value = -70.705725699067
$ws.Range("C20").Value = 2
$ws.Range("E20").Value = 100
$ws.Range("G20").Value = 2
$ws.Range("H20").Value = 200
$ws.Range("I20").Value = 67
$ws.Range("J20").Value = 46
$ws.Range("K20").Value = 45.652173913043
$ws.Range("L20").Value = 42.553191489361
$ws.Range("M20").Value = 116.129032258065
$ws.Range("N20").Value = -89.348171701112
$ws.Range("C21").Value = 15
$ws.Range("D21").Value = 23
$ws.Range("E21").Value = -34.782608695652
$ws.Range("F21").Value = 99
$ws.Range("G21").Value = 91
$ws.Range("H21").Value = 8.791208791208
$ws.Range("I21").Value = 1129
$ws.Range("J21").Value = 960
$ws.Range("K21").Value = 17.604166666666
$ws.Range("L21").Value = 18.967334035827
$ws.Range("M21").Value = 15.204081632653
$ws.Range("N21").Value = -76.874231872183
$ws.Range("D22").Value = "0"
$ws.Range("E22").Value = "***.*"
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 32
$ws.Range("K22").Value = 23.076923076923
$ws.Range("L22").Value = 52.380952380952
$ws.Range("M22").Value = -11.111111111111
$ws.Range("C24").Value = 12
$ws.Range("D24").Value = 12
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 65
$ws.Range("G24").Value = 79
$ws.Range("H24").Value = -17.721518987341
$ws.Range("I24").Value = 1210
$ws.Range("J24").Value = 989
$ws.Range("K24").Value = 22.345803842264
$ws.Range("L24").Value = -6.274206041828
$ws.Range("M24").Value = 95.476575121163
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = -66.666666666666
$ws.Range("F25").Value = 21
$ws.Range("H25").Value = 10.526315789473
$ws.Range("I25").Value = 236
$ws.Range("J25").Value = 181
$ws.Range("K25").Value = 30.386740331491
$ws.Range("L25").Value = 56.291390728476
$ws.Range("M25").Value = 5.357142857142
$ws.Range("C27").Value = "0"
$ws.Range("D27").Value = "0"
$ws.Range("E27").Value = "***.*"
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = -60
$ws.Range("D30").Value = 1
$ws.Range("E30").Value = -100
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = -100
$ws.Range("J30").Value = 7
$ws.Range("K30").Value = 71.428571428571
